$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I12").Value = "ba"
$ws.Range("J12").Value = "Appreciation"
$ws.Range("I14").Value = "sv"
$ws.Range("J14").Value = "Statement-opinion"
$ws.Range("I23").Value = "%"
$ws.Range("J23").Value = "Uninterpretable"
$ws.Range("I37").Value = "sv"
$ws.Range("J37").Value = "Statement-opinion"
$ws.Range("I45").Value = "aa"
$ws.Range("J45").Value = "Agree/Accept"
$ws.Range("I49").Value = "aa"
$ws.Range("J49").Value = "Agree/Accept"
$ws.Range("I53").Value = "sd"
$ws.Range("J53").Value = "Statement-non-opinion"
$ws.Range("I54").Value = "aa"
$ws.Range("J54").Value = "Agree/Accept"
$ws.Range("I56").Value = "sv"
$ws.Range("J56").Value = "Statement-opinion"
$ws.Range("I59").Value = "sd"
$ws.Range("J59").Value = "Statement-non-opinion"
$ws.Range("I60").Value = "sv"
$ws.Range("J60").Value = "Statement-opinion"
$ws.Range("I72").Value = "b"
$ws.Range("J72").Value = "Acknowledge (Backchannel)"
$ws.Range("I74").Value = "sv"
$ws.Range("J74").Value = "Statement-opinion"
$ws.Range("I75").Value = "sv"
$ws.Range("J75").Value = "Statement-opinion"
$ws.Range("I82").Value = "sv"
$ws.Range("J82").Value = "Statement-opinion"
$ws.Range("I85").Value = "sv"
$ws.Range("J85").Value = "Statement-opinion"
$ws.Range("I86").Value = "aa"
$ws.Range("J86").Value = "Agree/Accept"
$ws.Range("I88").Value = "sv"
$ws.Range("J88").Value = "Statement-opinion"
$ws.Range("I90").Value = "sd"
$ws.Range("J90").Value = "Statement-non-opinion"
$ws.Range("I97").Value = "b"
$ws.Range("J97").Value = "Acknowledge (Backchannel)"
$ws.Range("I103").Value = "sd"
$ws.Range("J103").Value = "Statement-non-opinion"
$ws.Range("I105").Value = "sd"
$ws.Range("J105").Value = "Statement-non-opinion"
$ws.Range("I123").Value = "aa"
$ws.Range("J123").Value = "Agree/Accept"
$ws.Range("I125").Value = "ba"
$ws.Range("J125").Value = "Appreciation"
$ws.Range("I138").Value = "sd"
$ws.Range("J138").Value = "Statement-non-opinion"
$ws.Range("I144").Value = "sd"
$ws.Range("J144").Value = "Statement-non-opinion"
$ws.Range("I150").Value = "aa"
$ws.Range("J150").Value = "Agree/Accept"
$ws.Range("I165").Value = "aa"
$ws.Range("J165").Value = "Agree/Accept"
$ws.Range("I179").Value = "sd"
$ws.Range("J179").Value = "Statement-non-opinion"
$ws.Range("I188").Value = "sd"
$ws.Range("J188").Value = "Statement-non-opinion"
$ws.Range("I189").Value = "sd"
$ws.Range("J189").Value = "Statement-non-opinion"
$ws.Range("I194").Value = "sv"
$ws.Range("J194").Value = "Statement-opinion"
$ws.Range("I204").Value = "sv"
$ws.Range("J204").Value = "Statement-opinion"
$ws.Range("I223").Value = "aa"
$ws.Range("J223").Value = "Agree/Accept"
$ws.Range("I226").Value = "sv"
$ws.Range("J226").Value = "Statement-opinion"
$ws.Range("I232").Value = "sv"
$ws.Range("J232").Value = "Statement-opinion"
$ws.Range("I237").Value = "sd"
$ws.Range("J237").Value = "Statement-non-opinion"
$ws.Range("I243").Value = "sv"
$ws.Range("J243").Value = "Statement-opinion"
$ws.Range("I244").Value = "sd"
$ws.Range("J244").Value = "Statement-non-opinion"
$ws.Range("I245").Value = "sd"
$ws.Range("J245").Value = "Statement-non-opinion"
$ws.Range("I252").Value = "b"
$ws.Range("J252").Value = "Acknowledge (Backchannel)"
$ws.Range("I266").Value = "sv"
$ws.Range("J266").Value = "Statement-opinion"
$ws.Range("I267").Value = "sd"
$ws.Range("J267").Value = "Statement-non-opinion"
$ws.Range("I272").Value = "sd"
$ws.Range("J272").Value = "Statement-non-opinion"
$ws.Range("I287").Value = "sv"
$ws.Range("J287").Value = "Statement-opinion"
$ws.Range("I294").Value = "sv"
$ws.Range("J294").Value = "Statement-opinion"
$ws.Range("I301").Value = "ba"
$ws.Range("J301").Value = "Appreciation"
$ws.Range("I307").Value = "%"
$ws.Range("J307").Value = "Uninterpretable"
$ws.Range("I316").Value = "sd"
$ws.Range("J316").Value = "Statement-non-opinion"
$ws.Range("I331").Value = "sd"
$ws.Range("J331").Value = "Statement-non-opinion"
$ws.Range("I332").Value = "sd"
$ws.Range("J332").Value = "Statement-non-opinion"
$ws.Range("I335").Value = "sd"
$ws.Range("J335").Value = "Statement-non-opinion"
$ws.Range("I337").Value = "b"
$ws.Range("J337").Value = "Acknowledge (Backchannel)"
$ws.Range("I344").Value = "sv"
$ws.Range("J344").Value = "Statement-opinion"
$ws.Range("I363").Value = "aa"
$ws.Range("J363").Value = "Agree/Accept"
$ws.Range("I381").Value = "aa"
$ws.Range("J381").Value = "Agree/Accept"
$ws.Range("I389").Value = "sv"
$ws.Range("J389").Value = "Statement-opinion"
$ws.Range("I391").Value = "sv"
$ws.Range("J391").Value = "Statement-opinion"
$ws.Range("I393").Value = "sd"
$ws.Range("J393").Value = "Statement-non-opinion"
$ws.Range("I394").Value = "b"
$ws.Range("J394").Value = "Acknowledge (Backchannel)"
$ws.Range("I401").Value = "aa"
$ws.Range("J401").Value = "Agree/Accept"
$ws.Range("I411").Value = "aa"
$ws.Range("J411").Value = "Agree/Accept"
$ws.Range("I412").Value = "aa"
$ws.Range("J412").Value = "Agree/Accept"
$ws.Range("I414").Value = "aa"
$ws.Range("J414").Value = "Agree/Accept"
$ws.Range("I418").Value = "sv"
$ws.Range("J418").Value = "Statement-opinion"
$ws.Range("I419").Value = "aa"
$ws.Range("J419").Value = "Agree/Accept"
$ws.Range("I424").Value = "ba"
$ws.Range("J424").Value = "Appreciation"
$ws.Range("I430").Value = "sv"
$ws.Range("J430").Value = "Statement-opinion"
$ws.Range("I438").Value = "b"
$ws.Range("J438").Value = "Acknowledge (Backchannel)"
$ws.Range("I466").Value = "sv"
$ws.Range("J466").Value = "Statement-opinion"
$ws.Range("I469").Value = "b"
$ws.Range("J469").Value = "Acknowledge (Backchannel)"
$ws.Range("I490").Value = "sd"
$ws.Range("J490").Value = "Statement-non-opinion"
$ws.Range("I491").Value = "sv"
$ws.Range("J491").Value = "Statement-opinion"
$ws.Range("I499").Value = "sv"
$ws.Range("J499").Value = "Statement-opinion"
$ws.Range("I501").Value = "sv"
$ws.Range("J501").Value = "Statement-opinion"
$ws.Range("I503").Value = "sv"
$ws.Range("J503").Value = "Statement-opinion"
$ws.Range("I508").Value = "sv"
$ws.Range("J508").Value = "Statement-opinion"
$ws.Range("I520").Value = "sd"
$ws.Range("J520").Value = "Statement-non-opinion"
$ws.Range("I523").Value = "aa"
$ws.Range("J523").Value = "Agree/Accept"
$ws.Range("I524").Value = "sd"
$ws.Range("J524").Value = "Statement-non-opinion"
$ws.Range("I525").Value = "sv"
$ws.Range("J525").Value = "Statement-opinion"
$ws.Range("I528").Value = "sd"
$ws.Range("J528").Value = "Statement-non-opinion"
$ws.Range("I532").Value = "ba"
$ws.Range("J532").Value = "Appreciation"
$ws.Range("I545").Value = "sd"
$ws.Range("J545").Value = "Statement-non-opinion"
$ws.Range("I548").Value = "aa"
$ws.Range("J548").Value = "Agree/Accept"
$ws.Range("I549").Value = "qy"
$ws.Range("J549").Value = "Yes-No-Question"
